$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "03:03 Our result class have " / "a " / "generic " / "model and a throwable"
# were split across four runs. Re-typing the whole (unchanged) sentence via
# Find/Replace collapses it back into a single run, matching the target.
$null = $d.Content.Find.Execute(
    "03:03 Our result class have a generic model and a throwable",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "03:03 Our result class have a generic model and a throwable", 2)

# --- Change 2 -------------------------------------------------------------
# "getPreviouslyCalculateddata" -> "getPreviousLoveCalculations", and the
# following run " from the repository" becomes two runs: " " and
# "from the repository".
$wordRng = $d.Content
$found = $wordRng.Find.Execute(
    "getPreviouslyCalculateddata",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "getPreviousLoveCalculations", 2)

# Locate " from the repository" right after the word we just replaced.
$tailRng = $d.Range($wordRng.End, $d.Content.End)
$null = $tailRng.Find.Execute(
    " from the repository",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)

# Split that run into " " and "from the repository" by nudging a formatting
# property on just the leading space character and reverting it - this
# forces the run to break at that boundary without altering the text.
$spaceRng = $d.Range($tailRng.Start, $tailRng.Start + 1)
$spaceRng.Font.Bold = 1
$spaceRng.Font.Bold = 0
